$wb = $excel.ActiveWorkbook

$wsRedOak = $wb.Worksheets.Item("Red Oak")
$wsYellowPoplar = $wb.Worksheets.Item("Yellow Poplar")
$wsWhiteAsh = $wb.Worksheets.Item("White Ash")

# ---------------------------------------------------------------------------
# Red Oak (sheet1): add column J "delta_E" = SQRT((D-G)^2+(E-H)^2+(F-I)^2)
# data rows 2..22 (row 23 is a trailing blank row, no formula there)
# ---------------------------------------------------------------------------
$wsRedOak.Range("J1").Value = "delta_E"
$wsRedOak.Range("J1").HorizontalAlignment = -4108
$wsRedOak.Range("J2").Formula = "=SQRT((D2-G2)^2+(E2-H2)^2+(F2-I2)^2)"
$wsRedOak.Range("J3:J22").Formula = "=SQRT((D3-G3)^2+(E3-H3)^2+(F3-I3)^2)"

# ---------------------------------------------------------------------------
# Yellow Poplar (sheet2): same new column, data rows 2..37.
# Row 23 was re-entered individually by the author (breaking the shared
# formula group there), so it is written separately from the rest.
# ---------------------------------------------------------------------------
$wsYellowPoplar.Range("J1").Value = "delta_E"
$wsYellowPoplar.Range("J1").HorizontalAlignment = -4108
$wsYellowPoplar.Range("J2").Formula = "=SQRT((D2-G2)^2+(E2-H2)^2+(F2-I2)^2)"
$wsYellowPoplar.Range("J3:J22").Formula = "=SQRT((D3-G3)^2+(E3-H3)^2+(F3-I3)^2)"
$wsYellowPoplar.Range("J23").Formula = "=SQRT((D23-G23)^2+(E23-H23)^2+(F23-I23)^2)"
$wsYellowPoplar.Range("J24:J37").Formula = "=SQRT((D24-G24)^2+(E24-H24)^2+(F24-I24)^2)"

# ---------------------------------------------------------------------------
# White Ash (sheet3): same new column, only data rows 2..3 have values.
# ---------------------------------------------------------------------------
$wsWhiteAsh.Range("J1").Value = "delta_E"
$wsWhiteAsh.Range("J1").HorizontalAlignment = -4108
$wsWhiteAsh.Range("J2").Formula = "=SQRT((D2-G2)^2+(E2-H2)^2+(F2-I2)^2)"
$wsWhiteAsh.Range("J3").Formula = "=SQRT((D3-G3)^2+(E3-H3)^2+(F3-I3)^2)"

# ---------------------------------------------------------------------------
# Restore each sheet's own cursor/selection (left over from working on it),
# then finish with "Red Oak" as the active tab / selection, matching the
# final view state recorded in the saved workbook.
# ---------------------------------------------------------------------------
$wsWhiteAsh.Activate()
$wsWhiteAsh.Range("G22").Select()

$wsYellowPoplar.Activate()
$wsYellowPoplar.Range("J1:J3").Select()

$wsRedOak.Activate()
$wsRedOak.Range("B2:F2").Select()
